# Swap the order of "System" and the email address in the "Recorded By"
# column (G) wherever both appear together as "System, dnasr281@gmail.com",
# turning it into "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# Column G is the "Recorded By" column; data rows run from row 2 downward.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
